$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1928
$ws1.Range("F3").Value = 526
$ws1.Range("F5").Value = 184
$ws1.Range("F6").Value = 2868
$ws1.Range("F7").Value = 200
$ws1.Range("F9").Value = 187
$ws1.Range("F10").Value = 1595
$ws1.Range("F11").Value = 575
$ws1.Range("F16").Value = 184
$ws1.Range("F19").Value = 229
$ws1.Range("F25").Value = 269
$ws1.Range("F26").Value = 81
$ws1.Range("F27").Value = 84
$ws1.Range("F28").Value = 4
$ws1.Range("F29").Value = 1859
$ws1.Range("F30").Value = 46
$ws1.Range("F31").Value = 436
$ws1.Range("F32").Value = 2
$ws1.Range("F33").Value = 115
$ws1.Range("F34").Value = 573
$ws1.Range("F36").Value = 320
$ws1.Range("F38").Value = 466

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1928
$ws4.Range("F4").Value = 526
$ws4.Range("F6").Value = 184
$ws4.Range("F7").Value = 2868
$ws4.Range("F8").Value = 200
$ws4.Range("F10").Value = 187
$ws4.Range("F11").Value = 1595
$ws4.Range("F12").Value = 575
$ws4.Range("F17").Value = 184
$ws4.Range("F20").Value = 229
$ws4.Range("F26").Value = 269
$ws4.Range("F27").Value = 81
$ws4.Range("F28").Value = 84
$ws4.Range("F29").Value = 4
$ws4.Range("F30").Value = 1859
$ws4.Range("F31").Value = 46
$ws4.Range("F32").Value = 436
$ws4.Range("F33").Value = 2
$ws4.Range("F34").Value = 115
$ws4.Range("F35").Value = 573
$ws4.Range("F37").Value = 320
$ws4.Range("F39").Value = 466
